$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update the two cell labels that were renamed in this revision
$ws.Range("B9").Value = "daily_buletim_id"
$ws.Range("B10").Value = "situation"

# Keep the selection where the author left it when saving
$ws.Range("B11").Select()
